$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.524.52"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "1.916.84"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.36"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4788"
$ws.Range("E7").Value = "  +2.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4100"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.85"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08049"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.39"
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "1.930.12"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.949"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.153"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.58"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001032"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.75"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "29.541.24"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.543"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.51"
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.206"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").Value = "2.131.00"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.54"
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.84"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.811"
$ws.Range("E29").Value = "  +6.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.137"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.94"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.055"
$ws.Range("E32").Value = "  +6.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09527"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.425"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.574"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.386"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06100"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02256"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.333"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.175"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5883"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.545"
$ws.Range("E42").Value = "  +7.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1848"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.14"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07978"
$ws.Range("E45").Value = "  +12.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.287"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.25"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5549"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.932"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.30"
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.74"
$ws.Range("E51").Value = "  -6.48%  "
